$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price column (D) to remain Text so numeric-looking values (e.g. "581.77")
# are not auto-converted to numbers by Excel, matching the original inline-string data.

# Update price (D) and volume/1h (E) columns for changed rows
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.502.39"
$ws.Range("E2").Value = "  +0.56%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.623.51"
$ws.Range("E3").Value = "  +1.33%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "581.77"
$ws.Range("E5").Value = "  +2.60%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.51"
$ws.Range("E6").Value = "  +2.15%  "
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("E8").Value = "  +0.25%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.53"
$ws.Range("E9").Value = "  -0.66%  "
$ws.Range("E10").Value = "  +0.88%  "
$ws.Range("E11").Value = "  +1.83%  "
$ws.Range("E12").Value = "  +3.24%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.098.09"
$ws.Range("E13").Value = "  +1.53%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.48"
$ws.Range("E14").Value = "  +14.10%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "60.488.23"
$ws.Range("E15").Value = "  +0.49%  "
$ws.Range("E16").Value = "  +1.66%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.628.68"
$ws.Range("E17").Value = "  +1.09%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.52"
$ws.Range("E18").Value = "  +1.66%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.72"
$ws.Range("E19").Value = "  +1.80%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "347.70"
$ws.Range("E20").Value = "  +0.65%  "
$ws.Range("E21").Value = "  -0.40%  "
$ws.Range("E22").Value = "  -0.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.530"
$ws.Range("E23").Value = "  -1.27%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.72"
$ws.Range("E24").Value = "  +1.18%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.999"
$ws.Range("E25").Value = "  +0.77%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.161"
$ws.Range("E26").Value = "  +1.73%  "
$ws.Range("E27").Value = "  +6.17%  "
$ws.Range("E28").Value = "  +14.36%  "
$ws.Range("E29").Value = "  +2.19%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.64"
$ws.Range("E30").Value = "  +5.32%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "169.77"
$ws.Range("E31").Value = "  +5.83%  "
$ws.Range("E32").Value = "  +0.11%  "
$ws.Range("E33").Value = "  +0.98%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.07"
$ws.Range("E34").Value = "  +12.35%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.46"
$ws.Range("E35").Value = "  +5.53%  "
$ws.Range("E36").Value = "  +8.89%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.66"
$ws.Range("E37").Value = "  +4.44%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "330.61"
$ws.Range("E38").Value = "  +12.64%  "
$ws.Range("E39").Value = "  +5.80%  "
$ws.Range("E40").Value = "  +3.05%  "
$ws.Range("E41").Value = "  +1.33%  "
$ws.Range("E42").Value = "  +6.62%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "133.08"
$ws.Range("E43").Value = "  -3.63%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0996"
$ws.Range("E44").Value = "  +2.03%  "

# Rows 45-47 reorder: InjectiveProtocol moves to 45, EnergySwap to 46, FirstDigitalUSD to 47
$ws.Range("B45").Value = "InjectiveProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "20.57"
$ws.Range("E45").Value = "  +4.81%  "

$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "20.07"
$ws.Range("E46").Value = "  +2.55%  "

$ws.Range("B47").Value = "FirstDigitalUSD"
$ws.Range("C47").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.999"
$ws.Range("E47").Value = "  +0.24%  "

# remaining rows after 47
$ws.Range("E48").Value = "  +2.45%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.609"
$ws.Range("E49").Value = "  +0.77%  "
$ws.Range("E50").Value = "  +2.17%  "
$ws.Range("E51").Value = "  +0.76%  "
